# Rename the inline picture shapes that live in the document's headers and
# footers:
#   - the two Pearson Edexcel logo pictures (descr = "...PearsonLogo.png")
#     go from name "image1.png" -> name "image2.png"
#   - the BTec logo picture (descr = "BTec_Logo-Orange") goes from
#     name "image2.jpg" -> name "image1.jpg"
#
# InlineShape has no settable .Name in the Word object model, so each
# picture is temporarily converted to a floating Shape (which does expose
# .Name), renamed, then converted back to an inline shape so the
# <wp:inline> layout is preserved.

$d = $word.ActiveDocument

function Rename-InlineShapeByDescr($range, $descr, $newName) {
    if ($null -eq $range) { return }
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $ish = $shapes.Item($i)
        if ($ish.AlternativeText -eq $descr) {
            $shp = $ish.ConvertToShape()
            $shp.Name = $newName
            $shp.ConvertToInlineShape() | Out-Null
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($t = 1; $t -le 3; $t++) {
        $ftr = $sec.Footers.Item($t)
        if ($ftr.Exists) {
            Rename-InlineShapeByDescr $ftr.Range `
                "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" `
                "image2.png"
        }

        $hdr = $sec.Headers.Item($t)
        if ($hdr.Exists) {
            Rename-InlineShapeByDescr $hdr.Range "BTec_Logo-Orange" "image1.jpg"
        }
    }
}
